$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-label the second column header from "@name" to "name"
$ws.Range("B1").Value = "name"

# Leave the selection where the author's session ended up
[void]$ws.Range("H13").Select()
